# Update column F (dSF) values for the affected rows to reflect the
# repulled / recalculated data, as described in the commit message
# "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -7
    3  = -8
    4  = 3
    5  = 1
    6  = -5
    7  = -3
    8  = -1
    14 = -1
    17 = -2
    22 = 0
    24 = -1
    26 = -14
    27 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
